# Auto-generated edit script: updates leve-profit figures across sheets
# per the scheduled-runner price refresh (see commit message).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 286.08334
$ws.Range("I8").Value = 146.57143
$ws.Range("K8").Value = 439.71429
$ws.Range("M8").Value = -300.71429
$ws.Range("H15").Value = 1203.8108
$ws.Range("I15").Value = 1203.8108
$ws.Range("K15").Value = 3611.4324
$ws.Range("M15").Value = -3442.4324
$ws.Range("H94").Value = 1972.75
$ws.Range("I94").Value = 1070.5
$ws.Range("J94").Value = 2875
$ws.Range("K94").Value = 1070.5
$ws.Range("L94").Value = 2875
$ws.Range("M94").Value = -619.5
$ws.Range("N94").Value = -3777
$ws.Range("H99").Value = 202180.6
$ws.Range("I99").Value = 599
$ws.Range("J99").Value = 252576
$ws.Range("K99").Value = 1797
$ws.Range("L99").Value = 757728
$ws.Range("M99").Value = -299
$ws.Range("N99").Value = -760724
$ws.Range("H100").Value = 3993.625
$ws.Range("I100").Value = 3309.8
$ws.Range("K100").Value = 3309.8
$ws.Range("M100").Value = -2768.8
$ws.Range("H101").Value = 795
$ws.Range("J101").Value = 717.5
$ws.Range("L101").Value = 2152.5
$ws.Range("N101").Value = -5396.5
$ws.Range("H106").Value = 2244
$ws.Range("I106").Value = 2325.3333
$ws.Range("K106").Value = 2325.3333
$ws.Range("M106").Value = -1694.3333
$ws.Range("H113").Value = 7190.6665
$ws.Range("I113").Value = 7190.6665
$ws.Range("K113").Value = 7190.6665
$ws.Range("M113").Value = -3936.6665
$ws.Range("H116").Value = 9363.272000000001
$ws.Range("I116").Value = 8856.714
$ws.Range("J116").Value = 10249.75
$ws.Range("K116").Value = 8856.714
$ws.Range("L116").Value = 10249.75
$ws.Range("M116").Value = -5414.714
$ws.Range("N116").Value = -17133.75
$ws.Range("H131").Value = 5713.1113
$ws.Range("I131").Value = 1461.3182
$ws.Range("J131").Value = 24421
$ws.Range("K131").Value = 4383.9546
$ws.Range("L131").Value = 73263
$ws.Range("M131").Value = 656.0454
$ws.Range("N131").Value = -83343
$ws.Range("H132").Value = 3864.04
$ws.Range("I132").Value = 3812.8086
$ws.Range("K132").Value = 11438.4258
$ws.Range("M132").Value = -8908.425799999999
$ws.Range("H133").Value = 49193.08
$ws.Range("J133").Value = 49193.08
$ws.Range("L133").Value = 49193.08
$ws.Range("N133").Value = -59313.08
$ws.Range("H138").Value = 4662.4375
$ws.Range("I138").Value = 1578.037
$ws.Range("J138").Value = 6233.736
$ws.Range("K138").Value = 4734.111
$ws.Range("L138").Value = 18701.208
$ws.Range("M138").Value = 405.8890000000001
$ws.Range("N138").Value = -28981.208

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 31789.8
$ws.Range("J44").Value = 31789.8
$ws.Range("L44").Value = 31789.8
$ws.Range("N44").Value = -32765.8
$ws.Range("H55").Value = 43666.668
$ws.Range("J55").Value = 43666.668
$ws.Range("L55").Value = 43666.668
$ws.Range("N55").Value = -44296.668
$ws.Range("H61").Value = 2813.025
$ws.Range("I61").Value = 2269.9333
$ws.Range("J61").Value = 4442.3
$ws.Range("K61").Value = 2269.9333
$ws.Range("L61").Value = 4442.3
$ws.Range("M61").Value = -2057.9333
$ws.Range("N61").Value = -4866.3
$ws.Range("H63").Value = 9124.875
$ws.Range("I63").Value = 3000
$ws.Range("J63").Value = 9999.857
$ws.Range("K63").Value = 3000
$ws.Range("L63").Value = 9999.857
$ws.Range("M63").Value = -2314
$ws.Range("N63").Value = -11371.857
$ws.Range("H66").Value = 9124.875
$ws.Range("I66").Value = 3000
$ws.Range("J66").Value = 9999.857
$ws.Range("K66").Value = 15000
$ws.Range("L66").Value = 49999.285
$ws.Range("M66").Value = -11568
$ws.Range("N66").Value = -56863.285
$ws.Range("H97").Value = 1376.2
$ws.Range("I97").Value = 1506.9445
$ws.Range("J97").Value = 199.5
$ws.Range("K97").Value = 1506.9445
$ws.Range("L97").Value = 199.5
$ws.Range("M97").Value = -1010.9445
$ws.Range("N97").Value = -1191.5
$ws.Range("H102").Value = 1947.6471
$ws.Range("I102").Value = 1606.6666
$ws.Range("K102").Value = 1606.6666
$ws.Range("M102").Value = 15.33339999999998
$ws.Range("H136").Value = 2813.025
$ws.Range("I136").Value = 2269.9333
$ws.Range("J136").Value = 4442.3
$ws.Range("K136").Value = 6809.7999
$ws.Range("L136").Value = 13326.9
$ws.Range("M136").Value = -4259.7999
$ws.Range("N136").Value = -18426.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 695.1579
$ws.Range("I94").Value = 713.5
$ws.Range("J94").Value = 643.8
$ws.Range("K94").Value = 713.5
$ws.Range("L94").Value = 643.8
$ws.Range("M94").Value = -262.5
$ws.Range("N94").Value = -1545.8
$ws.Range("H99").Value = 2666
$ws.Range("I99").Value = 2666
$ws.Range("K99").Value = 2666
$ws.Range("M99").Value = -1168
$ws.Range("H117").Value = 52000
$ws.Range("J117").Value = 52000
$ws.Range("L117").Value = 52000
$ws.Range("N117").Value = -61178

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 17886.5
$ws.Range("J28").Value = 17886.5
$ws.Range("L28").Value = 17886.5
$ws.Range("N28").Value = -18376.5
$ws.Range("H31").Value = 403155.06
$ws.Range("I31").Value = 1541459
$ws.Range("J31").Value = 11863.094
$ws.Range("K31").Value = 1541459
$ws.Range("L31").Value = 11863.094
$ws.Range("M31").Value = -1541164
$ws.Range("N31").Value = -12453.094
$ws.Range("H34").Value = 403155.06
$ws.Range("I34").Value = 1541459
$ws.Range("J34").Value = 11863.094
$ws.Range("K34").Value = 1541459
$ws.Range("L34").Value = 11863.094
$ws.Range("M34").Value = -1541257
$ws.Range("N34").Value = -12267.094
$ws.Range("H58").Value = 7416.2666
$ws.Range("I58").Value = 2091.0667
$ws.Range("J58").Value = 12741.467
$ws.Range("K58").Value = 2091.0667
$ws.Range("L58").Value = 12741.467
$ws.Range("M58").Value = -1888.0667
$ws.Range("N58").Value = -13147.467
$ws.Range("H106").Value = 19671
$ws.Range("J106").Value = 19671
$ws.Range("L106").Value = 19671
$ws.Range("N106").Value = -22195
$ws.Range("H122").Value = 24535.867
$ws.Range("I122").Value = 42611.25
$ws.Range("J122").Value = 3878.2856
$ws.Range("K122").Value = 127833.75
$ws.Range("L122").Value = 11634.8568
$ws.Range("M122").Value = -125383.75
$ws.Range("N122").Value = -16534.8568
$ws.Range("H134").Value = 291983.06
$ws.Range("I134").Value = 2828.4285
$ws.Range("K134").Value = 8485.2855
$ws.Range("M134").Value = -5950.2855
$ws.Range("H136").Value = 7416.2666
$ws.Range("I136").Value = 2091.0667
$ws.Range("J136").Value = 12741.467
$ws.Range("K136").Value = 6273.2001
$ws.Range("L136").Value = 38224.401
$ws.Range("M136").Value = -3723.2001
$ws.Range("N136").Value = -43324.401
$ws.Range("H141").Value = 65982.89999999999
$ws.Range("J141").Value = 65982.89999999999
$ws.Range("L141").Value = 65982.89999999999
$ws.Range("N141").Value = -76342.89999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2779092.5
$ws.Range("J113").Value = 14715.714
$ws.Range("L113").Value = 44147.142
$ws.Range("N113").Value = -48487.142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 28101.143
$ws.Range("I102").Value = 45427.125
$ws.Range("J102").Value = 4999.8335
$ws.Range("K102").Value = 45427.125
$ws.Range("L102").Value = 4999.8335
$ws.Range("M102").Value = -43805.125
$ws.Range("N102").Value = -8243.833500000001
$ws.Range("H126").Value = 40008036
$ws.Range("I126").Value = 76926800
$ws.Range("K126").Value = 230780400
$ws.Range("M126").Value = -230777930
$ws.Range("H132").Value = 74268.64
$ws.Range("I132").Value = 2486.0908
$ws.Range("K132").Value = 7458.2724
$ws.Range("M132").Value = -4928.2724

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1912.6818
$ws.Range("I16").Value = 1926.55
$ws.Range("K16").Value = 1926.55
$ws.Range("M16").Value = -1756.55
$ws.Range("H17").Value = 150
$ws.Range("J17").Value = 150
$ws.Range("L17").Value = 150
$ws.Range("N17").Value = -490
$ws.Range("H18").Value = 4500
$ws.Range("I18").Value = 4500
$ws.Range("K18").Value = 4500
$ws.Range("M18").Value = -4328
$ws.Range("H93").Value = 2988.238
$ws.Range("I93").Value = 2905.2666
$ws.Range("J93").Value = 3195.6667
$ws.Range("K93").Value = 2905.2666
$ws.Range("L93").Value = 3195.6667
$ws.Range("M93").Value = -1657.2666
$ws.Range("N93").Value = -5691.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1661.8334
$ws.Range("I100").Value = 1894.2
$ws.Range("K100").Value = 3788.4
$ws.Range("M100").Value = -3247.4
$ws.Range("H122").Value = 4697.0312
$ws.Range("I122").Value = 4146.125
$ws.Range("K122").Value = 12438.375
$ws.Range("M122").Value = -9988.375
$ws.Range("H132").Value = 27778.773
$ws.Range("I132").Value = 1742.2188
$ws.Range("J132").Value = 97209.586
$ws.Range("K132").Value = 5226.6564
$ws.Range("L132").Value = 291628.758
$ws.Range("M132").Value = -2696.6564
$ws.Range("N132").Value = -296688.758
$ws.Range("H136").Value = 306408.84
$ws.Range("I136").Value = 326091.94
$ws.Range("J136").Value = 238611.56
$ws.Range("K136").Value = 978275.8200000001
$ws.Range("L136").Value = 715834.6799999999
$ws.Range("M136").Value = -975725.8200000001
$ws.Range("N136").Value = -720934.6799999999
